$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 27
$ws.Range("B2").Value = 78
$ws.Range("B3").Value = 40
$ws.Range("B4").Value = 67
$ws.Range("B5").Value = 55
$ws.Range("B6").Value = 68
$ws.Range("B7").Value = 51
$ws.Range("B8").Value = 94
$ws.Range("B9").Value = 35
$ws.Range("B10").Value = 12
